$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New comment rows to append below the existing table (IDs 19-23), matching
# the existing table layout: column A = ID (number), column B = Comment
# (text), column C = Time (text date string, e.g. "04-10-2023").
$newRows = @(
    @(19, "let us test the comment system", "04-10-2023"),
    @(20, "test again", "04-10-2023"),
    @(21, "test it ", "04-10-2023"),
    @(22, "this ID should be 22", "04-10-2023"),
    @(23, "test 23", "04-10-2023")
)

$startRow = 21
$endRow = $startRow + $newRows.Length - 1

# Format column C as Text first so the date-like strings ("04-10-2023")
# are kept as literal text instead of being auto-converted to date serials,
# consistent with the existing rows above them.
$ws.Range("C$startRow`:C$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}
